# Rocky intertidal titrations 20240525
# Appends 8 new titration result rows (239-246) to the CRMAccuracyData sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Copy formatting from the last existing data row (238) down onto the new
#    rows so date formatting (column A) and the highlighted CRM-value style
#    (column C, only needed for the first two new rows) match the workbook's
#    existing conventions.
# ---------------------------------------------------------------------------
$ws.Range("A238").Copy() | Out-Null
$ws.Range("A239:A246").PasteSpecial(-4122) | Out-Null

$ws.Range("C238").Copy() | Out-Null
$ws.Range("C239:C240").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Column A - sample dates (stored as date serials)
# ---------------------------------------------------------------------------
$ws.Range("A239").Value = 45404
$ws.Range("A240").Value = 45434
$ws.Range("A241").Value = 45433
$ws.Range("A242").Value = 45435
$ws.Range("A243").Value = 45437
$ws.Range("A244").Value = 45437
$ws.Range("A245").Value = 45437
$ws.Range("A246").Value = 45437

# ---------------------------------------------------------------------------
# 3. Column B - measured CRM value
# ---------------------------------------------------------------------------
$ws.Range("B239").Value = 2231.97814827924
$ws.Range("B240").Value = 2228.69896965114
$ws.Range("B241").Value = 2216.9235939999999
$ws.Range("B242").Value = 2210.7065153475
$ws.Range("B243").Value = 2244.7447187303901
$ws.Range("B244").Value = 2232.7867714721801
$ws.Range("B245").Value = 2223.0452739453399
$ws.Range("B246").Value = 2228.5157841463201

# ---------------------------------------------------------------------------
# 4. Column C - accepted batch value
# ---------------------------------------------------------------------------
$ws.Range("C239").Value = 2230.52
$ws.Range("C240").Value = 2215.13
$ws.Range("C241").Value = 2215.13
$ws.Range("C242").Value = 2215.13
$ws.Range("C243").Value = 2215.13
$ws.Range("C244").Value = 2215.13
$ws.Range("C245").Value = 2215.13
$ws.Range("C246").Value = 2215.13

# ---------------------------------------------------------------------------
# 5. Column D - % off, a shared formula continuing the existing pattern
#    (100*(B-C)/C). Row 241 keeps a manually entered, hard-coded number in
#    the source workbook rather than a live formula, so it is written as a
#    plain value.
# ---------------------------------------------------------------------------
$ws.Range("D239:D240").Formula = "=100*(B239-C239)/C239"
$ws.Range("D241").Value = 0.080970159
$ws.Range("D242:D246").Formula = "=100*(B242-C242)/C242"

# ---------------------------------------------------------------------------
# 6. Column E - batch #
# ---------------------------------------------------------------------------
$ws.Range("E239").Value = 183
$ws.Range("E240").Value = 202
$ws.Range("E241").Value = 202
$ws.Range("E242").Value = 202
$ws.Range("E243").Value = 202
$ws.Range("E244").Value = 202
$ws.Range("E245").Value = 202
$ws.Range("E246").Value = 202

# ---------------------------------------------------------------------------
# 7. Column F - notes (new shared strings get appended automatically)
# ---------------------------------------------------------------------------
$ws.Range("F239").Value = "CRM opened 3/31/2024"
$ws.Range("F240").Value = "CRM opened 5/21/2024"
$ws.Range("F241").Value = "CRM opened 05/21/2024"
$ws.Range("F242").Value = "CRM opened 05/21/2024"
$ws.Range("F243").Value = "CRM opened 5/21 , opened new batch of acid"
$ws.Range("F244").Value = "CRM opened 05/21/2024"
$ws.Range("F245").Value = "CRM opened 05/21/2024"
$ws.Range("F246").Value = "CRM opened 05/21/2024"

# ---------------------------------------------------------------------------
# 8. Update the view so the selection tracks the newly added rows, mirroring
#    where Excel leaves the cursor after the edits were made.
# ---------------------------------------------------------------------------
$ws.Range("A241").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E251").Select() | Out-Null
$wb.Save()
